$d = $word.ActiveDocument

# Locate the "Expert methodology validated at highest judicial level" paragraph
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Expert methodology validated at highest judicial level*") {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -eq -1) {
    Write-Output "ERROR: anchor paragraph not found"
} else {
    $anchor = $d.Paragraphs.Item($anchorIndex)

    # Insert first new paragraph right after the anchor paragraph
    $anchor.Range.InsertParagraphAfter()
    $p1 = $d.Paragraphs.Item($anchorIndex + 1)
    $p1.Range.Text = "• Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions"

    # Insert second new paragraph right after the first new paragraph
    $p1.Range.InsertParagraphAfter()
    $p2 = $d.Paragraphs.Item($anchorIndex + 2)
    $p2.Range.Text = "• 178% accuracy improvement in racial classification algorithms"

    # Bold + color the "178%" portion of the second new paragraph
    $p2Range = $p2.Range
    $boldStart = $p2Range.Start + 2
    $boldEnd = $boldStart + 4
    $boldRange = $d.Range($boldStart, $boldEnd)
    $boldRange.Font.Bold = 1
    $boldRange.Font.Color = 5258796

    Write-Output "Inserted paragraphs after index $anchorIndex"
}
